$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 446.3243
$ws.Range("I19").Value = 387.5625
$ws.Range("J19").Value = 491.09525
$ws.Range("K19").Value = 387.5625
$ws.Range("L19").Value = 491.09525
$ws.Range("M19").Value = -212.5625
$ws.Range("N19").Value = -841.0952500000001
$ws.Range("H33").Value = 204.80952
$ws.Range("I33").Value = 148.33333
$ws.Range("K33").Value = 148.33333
$ws.Range("M33").Value = 80.66667000000001
$ws.Range("H40").Value = 3139.077
$ws.Range("I40").Value = 3812.625
$ws.Range("J40").Value = 2061.4
$ws.Range("K40").Value = 3812.625
$ws.Range("L40").Value = 2061.4
$ws.Range("M40").Value = -3637.625
$ws.Range("N40").Value = -2411.4
$ws.Range("H100").Value = 2252.7273
$ws.Range("I100").Value = 2016
$ws.Range("K100").Value = 2016
$ws.Range("M100").Value = -1475
$ws.Range("H131").Value = 10756.429
$ws.Range("I131").Value = 1765
$ws.Range("J131").Value = 17500
$ws.Range("K131").Value = 5295
$ws.Range("L131").Value = 52500
$ws.Range("M131").Value = -255
$ws.Range("N131").Value = -62580

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 30177.5
$ws.Range("J24").Value = 30177.5
$ws.Range("L24").Value = 30177.5
$ws.Range("N24").Value = -30925.5
$ws.Range("H32").Value = 2523.739
$ws.Range("I32").Value = 2026.641
$ws.Range("J32").Value = 5293.2856
$ws.Range("K32").Value = 2026.641
$ws.Range("L32").Value = 5293.2856
$ws.Range("M32").Value = -1739.641
$ws.Range("N32").Value = -5867.2856
$ws.Range("H100").Value = 30177.5
$ws.Range("J100").Value = 30177.5
$ws.Range("L100").Value = 30177.5
$ws.Range("N100").Value = -32341.5
$ws.Range("H132").Value = 3733.2273
$ws.Range("I132").Value = 3236.7273
$ws.Range("J132").Value = 4229.727
$ws.Range("K132").Value = 9710.1819
$ws.Range("L132").Value = 12689.181
$ws.Range("M132").Value = -7180.1819
$ws.Range("N132").Value = -17749.181
$ws.Range("H135").Value = 29500
$ws.Range("J135").Value = 29500
$ws.Range("L135").Value = 29500
$ws.Range("N135").Value = -39640

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3602.0715
$ws.Range("I134").Value = 1858.5714
$ws.Range("K134").Value = 5575.7142
$ws.Range("M134").Value = -3040.7142

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 168333740
$ws.Range("I4").Value = 625.25
$ws.Range("J4").Value = 505000000
$ws.Range("K4").Value = 625.25
$ws.Range("L4").Value = 505000000
$ws.Range("M4").Value = -513.25
$ws.Range("N4").Value = -505000224
$ws.Range("H31").Value = 2179
$ws.Range("I31").Value = 1469.2
$ws.Range("J31").Value = 2517
$ws.Range("K31").Value = 1469.2
$ws.Range("L31").Value = 2517
$ws.Range("M31").Value = -1174.2
$ws.Range("N31").Value = -3107
$ws.Range("H34").Value = 2179
$ws.Range("I34").Value = 1469.2
$ws.Range("J34").Value = 2517
$ws.Range("K34").Value = 1469.2
$ws.Range("L34").Value = 2517
$ws.Range("M34").Value = -1267.2
$ws.Range("N34").Value = -2921
$ws.Range("H58").Value = 956.05884
$ws.Range("I58").Value = 494.18918
$ws.Range("K58").Value = 494.18918
$ws.Range("M58").Value = -291.18918
$ws.Range("H132").Value = 2363.2432
$ws.Range("I132").Value = 1434.2609
$ws.Range("J132").Value = 3889.4285
$ws.Range("K132").Value = 4302.7827
$ws.Range("L132").Value = 11668.2855
$ws.Range("M132").Value = -1772.7827
$ws.Range("N132").Value = -16728.2855
$ws.Range("H136").Value = 956.05884
$ws.Range("I136").Value = 494.18918
$ws.Range("K136").Value = 1482.56754
$ws.Range("M136").Value = 1067.43246

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1669
$ws.Range("I60").Value = 170
$ws.Range("J60").Value = 5166.6665
$ws.Range("K60").Value = 510
$ws.Range("L60").Value = 15499.9995
$ws.Range("M60").Value = -259
$ws.Range("N60").Value = -16001.9995
$ws.Range("H131").Value = 2463.3206
$ws.Range("I131").Value = 355.55554
$ws.Range("J131").Value = 2738.2463
$ws.Range("K131").Value = 1066.66662
$ws.Range("L131").Value = 8214.7389
$ws.Range("M131").Value = 3973.33338
$ws.Range("N131").Value = -18294.7389

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2333.3333
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("H23").Value = 3983.3333
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 6966.6665
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 6966.6665
$ws.Range("M23").Value = -777
$ws.Range("N23").Value = -7412.6665
$ws.Range("H99").Value = 3092
$ws.Range("I99").Value = 3092
$ws.Range("K99").Value = 3092
$ws.Range("M99").Value = -846

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 50001
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H17").Value = 3000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H132").Value = 4023.2896
$ws.Range("I132").Value = 3653.0952
$ws.Range("J132").Value = 4480.5884
$ws.Range("K132").Value = 10959.2856
$ws.Range("L132").Value = 13441.7652
$ws.Range("M132").Value = -8429.285600000001
$ws.Range("N132").Value = -18501.7652
$ws.Range("H136").Value = 2619.0278
$ws.Range("I136").Value = 1949.1818
$ws.Range("K136").Value = 5847.5454
$ws.Range("M136").Value = -3297.5454

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 40026100
$ws.Range("I2").Value = 57150144
$ws.Range("K2").Value = 57150144
$ws.Range("M2").Value = -57150032
$ws.Range("H41").Value = 5034.2856
$ws.Range("J41").Value = 5034.2856
$ws.Range("L41").Value = 5034.2856
$ws.Range("N41").Value = -5814.2856
$ws.Range("H45").Value = 7300.8
$ws.Range("J45").Value = 7300.8
$ws.Range("L45").Value = 7300.8
$ws.Range("N45").Value = -8282.799999999999
$ws.Range("H132").Value = 16668889
$ws.Range("I132").Value = 25001832
$ws.Range("J132").Value = 3005.7
$ws.Range("K132").Value = 75005496
$ws.Range("L132").Value = 9017.099999999999
$ws.Range("M132").Value = -75002966
$ws.Range("N132").Value = -14077.1
$ws.Range("H136").Value = 10102542
$ws.Range("I136").Value = 14493504
$ws.Range("J136").Value = 3328.7
$ws.Range("K136").Value = 43480512
$ws.Range("L136").Value = 9986.099999999999
$ws.Range("M136").Value = -43477962
$ws.Range("N136").Value = -15086.1
